$wb = $excel.ActiveWorkbook

# --- Fix "mdb_addr" / "save_addr" auto-appointment bug ---
# Cells that should read "auto" but were left with a stray numeric value
# (or, in one case, a placeholder blank string) on the "os" and "device"
# sheets.

$wsOs = $wb.Worksheets.Item("os")
$wsOs.Range("F6").Value = "auto"
$wsOs.Range("J9").Value = "auto"
$wsOs.Range("J11").Value = "auto"
$wsOs.Range("F12").Value = "auto"

$wsDevice = $wb.Worksheets.Item("device")
$wsDevice.Range("F4").Value = "auto"
$wsDevice.Range("F5").Value = "auto"
$wsDevice.Range("F7").Value = "auto"

# --- Selection / active-sheet bookkeeping ---
# Move the active-cell selection on "device" off of K4, then switch the
# active workbook tab from "device" over to "os".
$wsDevice.Range("H15").Select()
$wsOs.Activate()
